$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.930.15'
$ws.Range('E2').Value = '  +1.60%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.691.01'
$ws.Range('E3').Value = '  -0.30%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.83%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '315.57'
$ws.Range('E5').Value = '  -0.52%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.007'
$ws.Range('E6').Value = '  +0.53%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3953'
$ws.Range('E7').Value = '  +1.19%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3987'
$ws.Range('E8').Value = '  -2.11%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.448'
$ws.Range('E9').Value = '  -2.72%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '52.57'
$ws.Range('E10').Value = '  -0.78%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.010'
$ws.Range('E11').Value = '  +1.04%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08709'
$ws.Range('E12').Value = '  -1.23%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '25.46'
$ws.Range('E13').Value = '  -3.67%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.353'
$ws.Range('E14').Value = '  -1.64%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.00001334'
$ws.Range('E15').Value = '  -2.12%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.844'
$ws.Range('E16').Value = '  -4.65%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '1.702.95'
$ws.Range('E17').Value = '  +0.15%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '94.59'
$ws.Range('E18').Value = '  -3.62%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.07188'
$ws.Range('E19').Value = '  -0.22%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '20.35'
$ws.Range('E20').Value = '  -1.07%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '7.158'
$ws.Range('E21').Value = '  -2.26%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.007'
$ws.Range('E22').Value = '  +0.46%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.16'
$ws.Range('E23').Value = '  -1.32%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.932.52'
$ws.Range('E24').Value = '  +1.59%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.379'
$ws.Range('E25').Value = '  +1.85%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.879'
$ws.Range('E26').Value = '  -4.27%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.19'
$ws.Range('E27').Value = '  +1.08%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.076'
$ws.Range('E28').Value = '  +8.00%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '161.62'
$ws.Range('E29').Value = '  -4.20%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '148.79'
$ws.Range('E30').Value = '  +2.76%  '

$ws.Range('B31').Value = 'WEMIXTOKEN'
$ws.Range('C31').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.640'
$ws.Range('E31').Value = '  +20.39%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.044'
$ws.Range('E32').Value = '  -4.68%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.892.18'
$ws.Range('E33').Value = '  +0.26%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08459'
$ws.Range('E34').Value = '  -3.74%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.03094'
$ws.Range('E35').Value = '  +0.37%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.024'
$ws.Range('E36').Value = '  -2.72%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.989'
$ws.Range('E37').Value = '  -3.27%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2836'
$ws.Range('E38').Value = '  +0.80%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.09659'
$ws.Range('E39').Value = '  +5.32%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '10.76'
$ws.Range('E40').Value = '  -1.57%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8059'
$ws.Range('E41').Value = '  +1.12%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '13.88'
$ws.Range('E42').Value = '  -2.33%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.465'
$ws.Range('E43').Value = '  -1.60%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '16.97'
$ws.Range('E44').Value = '  -2.19%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.616'
$ws.Range('E45').Value = '  -2.08%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.7231'
$ws.Range('E46').Value = '  -0.38%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.224'
$ws.Range('E47').Value = '  -0.98%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.08873'
$ws.Range('E48').Value = '  +8.05%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.379'
$ws.Range('E49').Value = '  -1.38%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.006'
$ws.Range('E50').Value = '  +0.41%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '138.11'
$ws.Range('E51').Value = '  -1.62%  '
